$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Docentes responsaveis:" value row (B13/C13, no A value) is removed;
# everything below shifts up by one row.
$ws.Rows.Item(13).Delete()

# After the shift, several (now-misaligned) data cells need their text
# updated to match the new content for this revision.
$ws.Range("B10").Value2 = "5817344 - Livia Melo Carneiro"
$ws.Range("C10").Value2 = "5817344 - Livia Melo Carneiro"

$ws.Range("B13").Value2 = "Semestral"
$ws.Range("C13").Value2 = "Semestral"

# "01/01/2012" would be auto-converted to a date serial by Value2, so copy
# the already-text "01/01/2012" cells (row 8) which keeps it a text string
# with the correct B/C column styles.
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))

$ws.Range("B18").Value2 = "5817344 - Livia Melo Carneiro"
$ws.Range("C18").Value2 = "5817344 - Livia Melo Carneiro"

$ws.Range("B19").Value2 = "O aluno será avaliado através de duas provas escritas P1 e P2."
$ws.Range("C19").Value2 = "O aluno será avaliado através de duas provas escritas P1 e P2."

$ws.Range("B20").Value2 = "A nota final NF será (P1 + P2)/2 ."
$ws.Range("C20").Value2 = "A nota final NF será (P1 + P2)/2 ."

$ws.Range("B21").Value2 = "Prova escrita sobre toda matéria. A média final MF será a média da nota final NF e da nota obtida na recuperação NR: MF = (NF + NR)/2 ."
$ws.Range("C21").Value2 = "Prova escrita sobre toda matéria. A média final MF será a média da nota final NF e da nota obtida na recuperação NR: MF = (NF + NR)/2 ."
